$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (this pushes the TOT row and the following
# blank rows down by one, matching how a new entry was added through the
# "Tabella1" table).
$ws.Rows("14").Insert()

# Copy the date formatting from the row above so the new date cell reuses
# the existing date style instead of creating a new one.
$ws.Range("H13").Copy()
$ws.Range("H14").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the new time-tracking entry.
$ws.Range("E14").Value = "Gianluca"
$ws.Range("F14").Value = "Interno"
$ws.Range("G14").Value = "CM"
$ws.Range("H14").Value = 43550   # 26/03/2019
$ws.Range("I14").Value = 30

# Fix up the TOT formula (now on row 15) so it includes the new row.
$ws.Range("I15").Formula = "=SUM(I2:I14)"

# Match the selection state recorded in the saved workbook.
$ws.Range("I15").Select()
